$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data cleanup: remove rows for players that no longer belong in this
# player_per_game_df export ("LaMelo Ball" and "LeBron James").
$playersToRemove = @("LaMelo Ball", "LeBron James")

# Find the last used row/col to know the current extent of the data.
$lastRow = $ws.Cells.SpecialCells(11).Row  # xlCellTypeLastCell = 11

foreach ($playerName in $playersToRemove) {
    for ($r = $lastRow; $r -ge 1; $r--) {
        $cellValue = $ws.Cells.Item($r, 1).Value()
        if ($cellValue -eq $playerName) {
            $ws.Rows.Item($r).Delete()
            break
        }
    }
}
